# Add two new version-history rows ("1.7" and "1.6") to the top of the
# data in Table1 on Sheet1, pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right under the header row (row 1), shifting
# the existing data rows (old row 2 onward) down by two.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Grow the table ("Table1") so it covers the two new rows as well.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D9"))

# Fill row 3 first (v1.6) and then row 2 (v1.7) so new shared-string
# entries land in the same order as the authored workbook.
$ws.Range("A3").Value = "1.6"
$ws.Range("B3").Value = "New calculations available in calculator"
$ws.Range("C3").Value = "Adam Mohd Taufik"
$ws.Range("D3").Value = (Get-Date -Year 2026 -Month 1 -Day 23 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("D3").NumberFormat = "m/d/yy"

$ws.Range("A2").Value = "1.7"
$ws.Range("B2").Value = "Introduced Logging Planner app"
$ws.Range("C2").Value = "Adam Mohd Taufik"
$ws.Range("D2").Value = (Get-Date -Year 2026 -Month 1 -Day 30 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("D2").NumberFormat = "m/d/yy"
